# Updated cryptos list on Thu Sep 21 23:55:43 UTC 2023 with GitHub Actions
# Refresh the Price (column D) and Volume(1h) (column E) figures for the
# crypto table on the active sheet, row by row, to match the latest
# scraped values. Price cells that look like a plain decimal number are
# forced back to text (matching the original inline-string storage)
# by toggling the number format to Text before the write and clearing
# the one-off format afterwards so no stray cell style is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.609.13"
$ws.Range("E2").Value = "  -1.98%  "
$ws.Range("D3").Value = "1.586.57"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.09"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.508"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.70%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -2.45%  "
$ws.Range("E9").Value = "  -1.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.49"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0833"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "1.809.32"
$ws.Range("E12").Value = "  -2.28%  "
$ws.Range("D13").Value = "1.607.01"
$ws.Range("E13").Value = "  -1.00%  "
$ws.Range("E14").Value = "  -2.96%  "
$ws.Range("E15").Value = "  -3.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.76"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "26.593.11"
$ws.Range("E18").Value = "  -2.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.21"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.12%  "
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("E21").Value = "  -3.29%  "
$ws.Range("E22").Value = "  -2.77%  "
$ws.Range("E23").Value = "  -2.21%  "
$ws.Range("E24").Value = "  -2.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.40"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.19"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.30%  "
$ws.Range("E28").Value = "  -2.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.27"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0506"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.17%  "
$ws.Range("E31").Value = "  -1.82%  "
$ws.Range("E32").Value = "  -4.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.685"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +23.60%  "
$ws.Range("E34").Value = "  -3.12%  "
$ws.Range("D35").Value = "1.308.71"
$ws.Range("E35").Value = "  -2.76%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.48"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.20%  "
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.824"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.38"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.64%  "
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("E43").Value = "  -2.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.60"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.47%  "
$ws.Range("D45").Value = "1.722.68"
$ws.Range("E45").Value = "  -2.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.35"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.60"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.837"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -9.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0504"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("E50").Value = "  -2.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.53"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.71%  "